# Swap the species-record data between row 12 and row 13 on the
# "Artfynd" sheet. Only the columns that actually differ between the
# two rows are touched: A, B, E, F, G, H, Q, R, AC.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "AC")

foreach ($col in $cols) {
    $cell12 = $ws.Range($col + "12")
    $cell13 = $ws.Range($col + "13")

    $val12 = $cell12.Value2
    $val13 = $cell13.Value2

    $cell12.Value2 = $val13
    $cell13.Value2 = $val12
}
